$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Subtitle 2 (id=3): move up slightly ---
$subtitle = $s.Shapes.Item(1)
$subtitle.Top = 461.4478

# --- Rectangle 3 (id=4): reposition/resize + rewrap text onto two lines at 16pt ---
$rect3 = $s.Shapes.Item(2)
$rect3.Left = 103.74575
$rect3.Top = 371.99425
$rect3.Width = 193.4515
$rect3.Height = 46.04528
$rect3.TextFrame.TextRange.Text = "Teaching Biodiversity for" + [char]13 + "FET Life Sciences"
$rect3.TextFrame.TextRange.Font.Size = 16

# --- Rectangle 8 (id=9): reposition/resize + bump font size to 16pt ---
$rect8 = $s.Shapes.Item(5)
$rect8.Top = 451.08245
$rect8.Height = 46.04528
$rect8.TextFrame.TextRange.Font.Size = 16

# --- Title 1 (id=2): move up ---
$title = $s.Shapes.Item(7)
$title.Left = 17.09174
$title.Top = 44.97111

# --- Group 55 (id=56): reposition ---
$group = $s.Shapes.Item(9)
$group.Left = 31.43851
$group.Top = 228.77166

# Export the updated deck as PDF alongside the pptx save.
$p.SaveAs("/tmp/work/coverslide.pdf", 32)
